$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'57.965.59"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +2.84%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.044.13"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +2.28%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.12%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'518.57"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +3.23%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'140.89"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +5.04%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.07%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.444"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +4.07%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'7.46"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +2.12%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.111"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +5.86%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.368"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +5.02%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'3.596.18"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +3.05%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  +2.20%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'26.66"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +6.23%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.0000169"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +13.45%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'58.045.47"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +2.94%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'6.19"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +8.85%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'3.057.15"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +2.64%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'12.99"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +5.11%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'8.07"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +3.97%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'335.34"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +3.33%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'5.78"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +1.60%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.998"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.15%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'0.500"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +6.37%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'65.14"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +5.04%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'0.168"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +3.68%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'0.0₃0952"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +7.21%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'0.989"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -1.09%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'6.88"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +5.91%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'7.51"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +10.88%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'1.83"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +5.40%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'1.22"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +3.14%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'21.03"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +3.40%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = "'Monero"
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = "'156.85"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.13%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('B35').Value = "'NEARProtocol"
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = "'4.75"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +7.11%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'5.92"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +7.01%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'1.29"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +1.48%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'25.17"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +9.76%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.0690"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +2.74%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'3.089.03"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +2.58%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'37.63"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +3.98%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'3.89"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +9.50%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'FirstDigitalUSD"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'1.00"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.19%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'Mantle"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'0.665"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +3.89%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'2.321.00"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +3.39%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'1.45"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +4.03%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'1.00"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +1.67%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'6.05"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +5.37%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'0.0240"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +2.37%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'19.66"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +3.94%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'1.85"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -4.03%  "
$ws.Range('E51').Style = 'Normal'
